# Apply the commit's changes:
# - Remove ("yes") values from column F ("selected_for_analysis") for a set of
#   rows where the "needs_binning" flag was a false positive / redundant.
# - Update the frozen-pane/selection view state (scrolled back to top,
#   selection moved to F25 instead of F66).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows whose F-column "yes" marker must be cleared
$rowsToClear = @(9, 10, 11, 12, 13, 14, 18, 19, 20, 21, 22, 35, 37, 60, 61, 62, 63, 64, 65)

foreach ($r in $rowsToClear) {
    $ws.Range("F$r").ClearContents()
}

# Update the selected cell / view (frozen pane top-left resets to A2 and the
# active selection moves to F25)
$ws.Range("F25").Select()
